# Apply corrected IFRS financial figures (error solve ifrs list)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = 52649
$ws.Range("E2").Value = 294
$ws.Range("F2").Value = 389
$ws.Range("G2").Value = 411
$ws.Range("H2").Value = 339
$ws.Range("I2").Value = 338
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 19267
$ws.Range("L2").Value = 12596
$ws.Range("M2").Value = 6670
$ws.Range("N2").Value = 6663
$ws.Range("O2").Value = 7
$ws.Range("P2").Value = 1116
$ws.Range("Q2").Value = -604
$ws.Range("R2").Value = 758
$ws.Range("S2").Value = 391
$ws.Range("T2").Value = 30
$ws.Range("U2").Value = -634
$ws.Range("V2").Value = 4595
$ws.Range("W2").Value = 0.56
$ws.Range("X2").Value = 0.64
$ws.Range("Y2").Value = 5.36
$ws.Range("Z2").Value = 1.86
$ws.Range("AA2").Value = 188.84
$ws.Range("AB2").Value = 262.97
$ws.Range("AC2").Value = 1512
$ws.Range("AD2").Value = 20.62
$ws.Range("AE2").Value = 29839
$ws.Range("AF2").Value = 1.04
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.6
$ws.Range("AI2").Value = 33.07
$ws.Range("AJ2").Value = 22329802

# Row 3
$ws.Range("D3").Value = 42619
$ws.Range("E3").Value = 234
$ws.Range("F3").Value = 234
$ws.Range("G3").Value = 384
$ws.Range("H3").Value = 1472
$ws.Range("I3").Value = 1472
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 15040
$ws.Range("L3").Value = 10140
$ws.Range("M3").Value = 4900
$ws.Range("N3").Value = 4892
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 661
$ws.Range("Q3").Value = 748
$ws.Range("R3").Value = -170
$ws.Range("S3").Value = -2936
$ws.Range("T3").Value = 15
$ws.Range("U3").Value = 733
$ws.Range("V3").Value = 3036
$ws.Range("W3").Value = 0.55
$ws.Range("X3").Value = 3.45
$ws.Range("Y3").Value = 25.47
$ws.Range("Z3").Value = 8.58
$ws.Range("AA3").Value = 206.95
$ws.Range("AB3").Value = 648.65
$ws.Range("AC3").Value = 7336
$ws.Range("AD3").Value = 3.9
$ws.Range("AE3").Value = 36990
$ws.Range("AF3").Value = 0.77
$ws.Range("AG3").Value = 750
$ws.Range("AH3").Value = 2.62
$ws.Range("AI3").Value = 6.74
$ws.Range("AJ3").Value = 13228966

# Row 4
$ws.Range("D4").Value = 35588
$ws.Range("E4").Value = 305
$ws.Range("F4").Value = 305
$ws.Range("G4").Value = 214
$ws.Range("H4").Value = 81
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 15142
$ws.Range("L4").Value = 10812
$ws.Range("M4").Value = 4330
$ws.Range("N4").Value = 4320
$ws.Range("O4").Value = 9
$ws.Range("P4").Value = 661
$ws.Range("Q4").Value = 470
$ws.Range("R4").Value = -297
$ws.Range("S4").Value = 1036
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 467
$ws.Range("V4").Value = 4309
$ws.Range("W4").Value = 0.86
$ws.Range("X4").Value = 0.23
$ws.Range("Y4").Value = 1.73
$ws.Range("Z4").Value = 0.53
$ws.Range("AA4").Value = 249.74
$ws.Range("AB4").Value = 644.28
$ws.Range("AC4").Value = 602
$ws.Range("AD4").Value = 36.4
$ws.Range("AE4").Value = 33534
$ws.Range("AF4").Value = 0.65
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 2.28
$ws.Range("AI4").Value = 81
$ws.Range("AJ4").Value = 13228966

# Row 5
$ws.Range("D5").Value = 43060
$ws.Range("E5").Value = 340
$ws.Range("F5").Value = 340
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 251
$ws.Range("I5").Value = 250
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 12589
$ws.Range("L5").Value = 8961
$ws.Range("M5").Value = 3628
$ws.Range("N5").Value = 3618
$ws.Range("O5").Value = 10
$ws.Range("P5").Value = 661
$ws.Range("Q5").Value = 485
$ws.Range("R5").Value = 179
$ws.Range("S5").Value = -945
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 482
$ws.Range("V5").Value = 3420
$ws.Range("W5").Value = 0.79
$ws.Range("X5").Value = 0.58
$ws.Range("Y5").Value = 6.29
$ws.Range("Z5").Value = 1.81
$ws.Range("AA5").Value = 247
$ws.Range("AB5").Value = 669.64
$ws.Range("AC5").Value = 1888
$ws.Range("AD5").Value = 10.33
$ws.Range("AE5").Value = 28579
$ws.Range("AF5").Value = 0.68
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 3.08
$ws.Range("AI5").Value = 30.41
$ws.Range("AJ5").Value = 13228966

# Row 6
$ws.Range("D6").Value = 47140
$ws.Range("E6").Value = 505
$ws.Range("F6").Value = 505
$ws.Range("G6").Value = -38
$ws.Range("H6").Value = -87
$ws.Range("I6").Value = -88
$ws.Range("K6").Value = 10982
$ws.Range("L6").Value = 7922
$ws.Range("M6").Value = 3059
$ws.Range("N6").Value = 3048
$ws.Range("P6").Value = 661
$ws.Range("Q6").Value = -418
$ws.Range("R6").Value = -281
$ws.Range("S6").Value = -81
$ws.Range("T6").Value = 79
$ws.Range("U6").Value = -498
$ws.Range("V6").Value = 3434
$ws.Range("W6").Value = 1.07
$ws.Range("X6").Value = -0.18
$ws.Range("Y6").Value = -2.64
$ws.Range("Z6").Value = -0.73
$ws.Range("AA6").Value = 258.97
$ws.Range("AB6").Value = 643.65
$ws.Range("AC6").Value = -666
$ws.Range("AD6").Value = -42.18
$ws.Range("AE6").Value = 24076
$ws.Range("AF6").Value = 1.17
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 2.14
$ws.Range("AI6").Value = -86.18
$ws.Range("AJ6").Value = 13228966

# Row 7
$ws.Range("D7").Value = 44717
$ws.Range("E7").Value = 498
$ws.Range("G7").Value = 308
$ws.Range("H7").Value = 244
$ws.Range("I7").Value = 245
$ws.Range("K7").Value = 11446
$ws.Range("L7").Value = 8237
$ws.Range("M7").Value = 3209
$ws.Range("N7").Value = 3196
$ws.Range("P7").Value = 660
$ws.Range("Q7").Value = 788
$ws.Range("R7").Value = -96
$ws.Range("S7").Value = 198
$ws.Range("T7").Value = 82
$ws.Range("U7").Value = 294
$ws.Range("W7").Value = 1.11
$ws.Range("X7").Value = 0.55
$ws.Range("Y7").Value = 7.84
$ws.Range("Z7").Value = 2.18
$ws.Range("AA7").Value = 256.68
$ws.Range("AC7").Value = 1850
$ws.Range("AD7").Value = 9.51
$ws.Range("AE7").Value = 25249
$ws.Range("AF7").Value = 0.7
$ws.Range("AG7").Value = 600
$ws.Range("AH7").Value = 3.41
$ws.Range("AI7").Value = 32.43

# Row 8
$ws.Range("D8").Value = 45120
$ws.Range("E8").Value = 530
$ws.Range("G8").Value = 320
$ws.Range("H8").Value = 240
$ws.Range("I8").Value = 240
$ws.Range("K8").Value = 11010
$ws.Range("L8").Value = 7630
$ws.Range("M8").Value = 3370
$ws.Range("N8").Value = 3360
$ws.Range("P8").Value = 660
$ws.Range("Q8").Value = 280
$ws.Range("R8").Value = -740
$ws.Range("S8").Value = -80
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 1.18
$ws.Range("X8").Value = 0.53
$ws.Range("Y8").Value = 7.33
$ws.Range("Z8").Value = 2.19
$ws.Range("AA8").Value = 226.41
$ws.Range("AC8").Value = 1814
$ws.Range("AD8").Value = 8.76
$ws.Range("AE8").Value = 26542
$ws.Range("AF8").Value = 0.6
$ws.Range("AG8").Value = 600
$ws.Range("AH8").Value = 3.77
$ws.Range("AI8").Value = 33.07

# Row 9
$ws.Range("D9").Value = 46640
$ws.Range("E9").Value = 550
$ws.Range("G9").Value = 350
$ws.Range("H9").Value = 260
$ws.Range("I9").Value = 260
$ws.Range("K9").Value = 11240
$ws.Range("L9").Value = 7680
$ws.Range("M9").Value = 3560
$ws.Range("N9").Value = 3550
$ws.Range("P9").Value = 660
$ws.Range("Q9").Value = 610
$ws.Range("R9").Value = -740
$ws.Range("S9").Value = -80
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 1.18
$ws.Range("X9").Value = 0.56
$ws.Range("Y9").Value = 7.52
$ws.Range("Z9").Value = 2.34
$ws.Range("AA9").Value = 215.73
$ws.Range("AC9").Value = 1965
$ws.Range("AD9").Value = 8.09
$ws.Range("AE9").Value = 28043
$ws.Range("AF9").Value = 0.57
$ws.Range("AG9").Value = 600
$ws.Range("AH9").Value = 3.77
$ws.Range("AI9").Value = 30.53

